$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - add the missing End Time (C8); formulas in D8/E8/F8 recalc automatically
$ws.Cells.Item(8, 3).Value = 0.44791666666666669

# Row 9 - add Start Time (B9) and End Time (C9)
$ws.Cells.Item(9, 2).Value = 0.61805555555555558
$ws.Cells.Item(9, 3).Value = 0.71527777777777779

# Row 10 - add Date (A10), Start Time (B10), End Time (C10)
$ws.Cells.Item(10, 1).Value = 43338
$ws.Cells.Item(10, 2).Value = 0.78819444444444453
$ws.Cells.Item(10, 3).Value = 0.4993055555555555

# Row 11 - add Date (A11), Start Time (B11), End Time (C11)
$ws.Cells.Item(11, 1).Value = 43339
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0.2986111111111111

# Move the active selection to A12
$ws.Range("A12").Select()
